$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "('Canopy Vista', ['Land — Forest Plains', '({T}: Add {G} or {W}.)', 'Canopy Vista enters the battlefield tapped unless you control two or more basic lands.'])"
$ws.Range("A3").Value = "('Cinder Glade', ['Land — Mountain Forest', '({T}: Add {R} or {G}.)', 'Cinder Glade enters the battlefield tapped unless you control two or more basic lands.'])"
$ws.Range("A4").Value = "('Prairie Stream', ['Land — Plains Island', '({T}: Add {W} or {U}.)', 'Prairie Stream enters the battlefield tapped unless you control two or more basic lands.'])"
$ws.Range("A5").Value = "('Smoldering Marsh', ['Land — Swamp Mountain', '({T}: Add {B} or {R}.)', 'Smoldering Marsh enters the battlefield tapped unless you control two or more basic lands.'])"
$ws.Range("A6").Value = "('Sunken Hollow', ['Land — Island Swamp', '({T}: Add {U} or {B}.)', 'Sunken Hollow enters the battlefield tapped unless you control two or more basic lands.'])"

$ws.Rows("7:21").Delete()
